# Apply the requested edit: remove the "Fornecedor ID" column (column I)
# from the "Produtos" worksheet, shifting every subsequent column one
# position to the left (J->I, K->J, L->K, ... S->R).
#
# Deleting the column leaves the two data-validation-style cell comments
# anchored to their original (now stale) cells, so they are explicitly
# re-anchored one column to the left to track the header cells they
# describe ("Unidade Medida" and "Origem Mercadoria").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column I ("Fornecedor ID"); the remaining columns
# (J, K, L, ...) shift left into its place. Comment anchors, however,
# stay glued to their original cell addresses (K1, P1) instead of
# following the shift, so they now point one column too far right.
$ws.Columns("I:I").Delete()

# Collect the now-stale comments (still sitting at their pre-delete cell
# addresses) before touching them.
$comments = @()
foreach ($c in $ws.Comments) {
    $addr = $c.Parent.Address($false, $false)
    $comments += ,@($addr, $c.Text())
}

# Re-create each comment one column to the left of its stale position,
# so it once again matches the header cell it documents (J1 "Unidade
# Medida" and O1 "Origem Mercadoria").
foreach ($item in $comments) {
    $oldCell = $ws.Range($item[0])
    $text = $item[1]
    $newCell = $ws.Cells.Item($oldCell.Row, $oldCell.Column - 1)
    if ($newCell.Comment -ne $null) {
        $newCell.Comment.Delete()
    }
    $oldCell.Comment.Delete()
    $newCell.AddComment($text) | Out-Null
}

# Mirror the resulting selection: the whole column that now sits where
# the deleted column used to be.
$ws.Range("I1:I1048576").Select()
